# Insert a new weekly data row for "Ciboulette" (Vega Central Mapocho de
# Santiago) at row 245, pushing the existing rows 245:313 down to 246:314.
# This matches how the source data set is appended to week-over-week:
# newest record inserted in its sorted slot, later rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("245:245").Insert()

$ws.Range("A245").Value = 9
$ws.Range("B245").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C245").Value = 'Metropolitana'
$ws.Range("D245").Value = 44551
$ws.Range("E245").Value = 13
$ws.Range("F245").Value = 100112039
$ws.Range("G245").Value = 'Ciboulette'
$ws.Range("H245").Value = 'Sin especificar'
$ws.Range("I245").Value = 'Primera'
$ws.Range("J245").Value = 160
$ws.Range("K245").Value = 800
$ws.Range("L245").Value = 1000
$ws.Range("M245").Value = 900
$ws.Range("N245").Value = '$/docena de atados'
$ws.Range("O245").Value = 'Región Metropolitana'
$ws.Range("P245").Value = 300
$ws.Range("Q245").Value = 3
$ws.Range("R245").Value = 'Hortaliza'
